$d = $word.ActiveDocument

# Prevent Word's "smart quotes" autoformatting from mangling straight quotes
# in the text we insert.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (first) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$boldText = "Meta description"
$restText = ": Book of Dead is an Egyptian-themed online slot game with a high RTP and attractive maximum payout. Play it for free now and enjoy the gambling feature."

$insertPoint = $d.Range($metaStart, $metaStart)
$insertPoint.InsertAfter($boldText + $restText)

$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)
$boldRange.Bold = 1

# Leave a genuinely empty leading run, matching the style used throughout
# the rest of the document.
$emptyPoint = $d.Range($metaStart, $metaStart)
$emptyPoint.InsertBefore("")

# ------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document, and replace the italic meta-description paragraph's
#    text with the new image-prompt text.
# ------------------------------------------------------------------
$oldTitleText = "Play Book of Dead for Free - RTP, Maximum Payout & Gambling Feature"
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$dupTitleText = $dupTitlePara.Range.Text.TrimEnd([char]13)

if ($dupTitleText -eq $oldTitleText) {
    $dupTitlePara.Range.Delete()
}

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End

$newImagePrompt = "Create a cartoon-style feature image for the game `"Book of Dead`" featuring a happy Maya warrior with glasses. The image should have vibrant colors and showcase the ancient Egyptian theme of the game. The Maya warrior should be holding the Book of Dead and standing in front of the pyramids. The background should have a sunset hue with Egyptian hieroglyphics in the sky. Make sure to highlight the warrior's glasses which should be oversized and reflective. The overall style of the image should be fun and eye-catching, inviting players to dive into the adventure-packed world of `"Book of Dead.`""

$lastTextRange = $d.Range($lastStart, $lastEnd)
$lastTextRange.Text = $newImagePrompt
